# Insert a new weekly price observation as row 24 (pushing existing rows down).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new blank row at position 24 so existing data (rows 24-55) shifts to 25-56.
$ws.Rows.Item(24).Insert()

# Populate the new row with the new observation.
$ws.Cells.Item(24, 1).Value = 6
$ws.Cells.Item(24, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(24, 3).Value = "Metropolitana"
$ws.Cells.Item(24, 4).Value = 44973
$ws.Cells.Item(24, 5).Value = 13
$ws.Cells.Item(24, 6).Value = "Fruta"
$ws.Cells.Item(24, 7).Value = 100102
$ws.Cells.Item(24, 8).Value = "Cítricos"
$ws.Cells.Item(24, 9).Value = 100102006
$ws.Cells.Item(24, 10).Value = "Pomelo"
$ws.Cells.Item(24, 11).Value = "Start Ruby"
$ws.Cells.Item(24, 12).Value = "Segunda"
$ws.Cells.Item(24, 13).Value = 140
$ws.Cells.Item(24, 14).Value = 4000
$ws.Cells.Item(24, 15).Value = 5000
$ws.Cells.Item(24, 16).Value = 4500
$ws.Cells.Item(24, 17).Value = "$/caja 14 kilos"
$ws.Cells.Item(24, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(24, 19).Value = 321
$ws.Cells.Item(24, 20).Value = 14

# Apply the same date number format used by the rest of the "Fecha" column.
$ws.Cells.Item(24, 4).NumberFormat = $ws.Cells.Item(25, 4).NumberFormat
